$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2469.125
$ws.Range("I51").Value = 2380.2
$ws.Range("J51").Value = 2617.3333
$ws.Range("K51").Value = 2380.2
$ws.Range("L51").Value = 2617.3333
$ws.Range("M51").Value = -1896.2
$ws.Range("N51").Value = -3585.3333
$ws.Range("H98").Value = 1537.1765
$ws.Range("I98").Value = 1392
$ws.Range("J98").Value = 1744.5714
$ws.Range("K98").Value = 1392
$ws.Range("L98").Value = 1744.5714
$ws.Range("M98").Value = 106
$ws.Range("N98").Value = -4740.5714
$ws.Range("H111").Value = 111885
$ws.Range("I111").Value = 652.6667
$ws.Range("J111").Value = 167501.17
$ws.Range("K111").Value = 1958.0001
$ws.Range("L111").Value = 502503.51
$ws.Range("M111").Value = 1108.9999
$ws.Range("N111").Value = -508637.51
$ws.Range("H122").Value = 1537.1765
$ws.Range("I122").Value = 1392
$ws.Range("J122").Value = 1744.5714
$ws.Range("K122").Value = 4176
$ws.Range("L122").Value = 5233.7142
$ws.Range("M122").Value = -1726
$ws.Range("N122").Value = -10133.7142
$ws.Range("H129").Value = 1139.9608
$ws.Range("I129").Value = 566.6667
$ws.Range("J129").Value = 1175.7916
$ws.Range("K129").Value = 1700.0001
$ws.Range("L129").Value = 3527.3748
$ws.Range("M129").Value = 3299.9999
$ws.Range("N129").Value = -13527.3748
$ws.Range("H132").Value = 2123.7646
$ws.Range("I132").Value = 1917.7241
$ws.Range("J132").Value = 3318.8
$ws.Range("K132").Value = 5753.1723
$ws.Range("L132").Value = 9956.400000000001
$ws.Range("M132").Value = -3223.1723
$ws.Range("N132").Value = -15016.4
$ws.Range("H138").Value = 2176.8193
$ws.Range("J138").Value = 2552.2856
$ws.Range("L138").Value = 7656.8568
$ws.Range("N138").Value = -17936.8568
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1009.2143
$ws.Range("I2").Value = 956.3
$ws.Range("J2").Value = 1141.5
$ws.Range("K2").Value = 956.3
$ws.Range("L2").Value = 1141.5
$ws.Range("M2").Value = -843.3
$ws.Range("N2").Value = -1367.5
$ws.Range("H32").Value = 3899.75
$ws.Range("I32").Value = 3323.2366
$ws.Range("J32").Value = 11559.143
$ws.Range("K32").Value = 3323.2366
$ws.Range("L32").Value = 11559.143
$ws.Range("M32").Value = -3036.2366
$ws.Range("N32").Value = -12133.143
$ws.Range("H45").Value = 3913
$ws.Range("I45").Value = 5045.591
$ws.Range("K45").Value = 5045.591
$ws.Range("M45").Value = -4668.591
$ws.Range("H61").Value = 339810.28
$ws.Range("I61").Value = 8558.529
$ws.Range("J61").Value = 772985.6
$ws.Range("K61").Value = 8558.529
$ws.Range("L61").Value = 772985.6
$ws.Range("M61").Value = -8346.529
$ws.Range("N61").Value = -773409.6
$ws.Range("H116").Value = 1009.2143
$ws.Range("I116").Value = 956.3
$ws.Range("J116").Value = 1141.5
$ws.Range("K116").Value = 956.3
$ws.Range("L116").Value = 1141.5
$ws.Range("M116").Value = 1337.7
$ws.Range("N116").Value = -5729.5
$ws.Range("H132").Value = 2130511.5
$ws.Range("I132").Value = 1905.7632
$ws.Range("J132").Value = 11117959
$ws.Range("K132").Value = 5717.2896
$ws.Range("L132").Value = 33353877
$ws.Range("M132").Value = -3187.2896
$ws.Range("N132").Value = -33358937
$ws.Range("H133").Value = 44072.2
$ws.Range("J133").Value = 44072.2
$ws.Range("L133").Value = 44072.2
$ws.Range("N133").Value = -49132.2
$ws.Range("H135").Value = 51133.453
$ws.Range("J135").Value = 51133.453
$ws.Range("L135").Value = 51133.453
$ws.Range("N135").Value = -61273.453
$ws.Range("H136").Value = 339810.28
$ws.Range("I136").Value = 8558.529
$ws.Range("J136").Value = 772985.6
$ws.Range("K136").Value = 25675.587
$ws.Range("L136").Value = 2318956.8
$ws.Range("M136").Value = -23125.587
$ws.Range("N136").Value = -2324056.8
$ws.Range("H139").Value = 44502.5
$ws.Range("J139").Value = 44502.5
$ws.Range("L139").Value = 44502.5
$ws.Range("N139").Value = -54782.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1009.2143
$ws.Range("I3").Value = 956.3
$ws.Range("J3").Value = 1141.5
$ws.Range("K3").Value = 956.3
$ws.Range("L3").Value = 1141.5
$ws.Range("M3").Value = -842.3
$ws.Range("N3").Value = -1369.5
$ws.Range("H86").Value = 1729
$ws.Range("I86").Value = 1559.6666
$ws.Range("J86").Value = 2575.6667
$ws.Range("K86").Value = 1559.6666
$ws.Range("L86").Value = 2575.6667
$ws.Range("M86").Value = -436.6666
$ws.Range("N86").Value = -4821.6667
$ws.Range("H89").Value = 1729
$ws.Range("I89").Value = 1559.6666
$ws.Range("J89").Value = 2575.6667
$ws.Range("K89").Value = 7798.333000000001
$ws.Range("L89").Value = 12878.3335
$ws.Range("M89").Value = -2182.333000000001
$ws.Range("N89").Value = -24110.3335
$ws.Range("H134").Value = 17658.377
$ws.Range("I134").Value = 3327.2307
$ws.Range("J134").Value = 61494.824
$ws.Range("K134").Value = 9981.6921
$ws.Range("L134").Value = 184484.472
$ws.Range("M134").Value = -7446.6921
$ws.Range("N134").Value = -189554.472
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 196020.89
$ws.Range("I31").Value = 1530.5272
$ws.Range("J31").Value = 623899.7
$ws.Range("K31").Value = 1530.5272
$ws.Range("L31").Value = 623899.7
$ws.Range("M31").Value = -1235.5272
$ws.Range("N31").Value = -624489.7
$ws.Range("H34").Value = 196020.89
$ws.Range("I34").Value = 1530.5272
$ws.Range("J34").Value = 623899.7
$ws.Range("K34").Value = 1530.5272
$ws.Range("L34").Value = 623899.7
$ws.Range("M34").Value = -1328.5272
$ws.Range("N34").Value = -624303.7
$ws.Range("H133").Value = 36875.332
$ws.Range("J133").Value = 36875.332
$ws.Range("L133").Value = 36875.332
$ws.Range("N133").Value = -41935.332
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 555630.7
$ws.Range("I8").Value = 555630.7
$ws.Range("K8").Value = 1666892.1
$ws.Range("M8").Value = -1666753.1
$ws.Range("H107").Value = 458.09525
$ws.Range("I107").Value = 430.92856
$ws.Range("J107").Value = 512.4286
$ws.Range("K107").Value = 1292.78568
$ws.Range("L107").Value = 1537.2858
$ws.Range("M107").Value = 627.21432
$ws.Range("N107").Value = -5377.2858
$ws.Range("H131").Value = 2778734.5
$ws.Range("I131").Value = 7692971.5
$ws.Range("J131").Value = 1122.3914
$ws.Range("K131").Value = 23078914.5
$ws.Range("L131").Value = 3367.1742
$ws.Range("M131").Value = -23073874.5
$ws.Range("N131").Value = -13447.1742
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4412.8037
$ws.Range("I132").Value = 4773.575
$ws.Range("J132").Value = 3510.875
$ws.Range("K132").Value = 14320.725
$ws.Range("L132").Value = 10532.625
$ws.Range("M132").Value = -11790.725
$ws.Range("N132").Value = -15592.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -827
$ws.Range("H61").Value = 1235.1818
$ws.Range("I61").Value = 1133.7
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 1133.7
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -931.7
$ws.Range("N61").Value = -2654
$ws.Range("H113").Value = 1235.1818
$ws.Range("I113").Value = 1133.7
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1133.7
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 1036.3
$ws.Range("N113").Value = -6590
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 200000900
$ws.Range("I107").Value = 500000640
$ws.Range("J107").Value = 1066.3334
$ws.Range("K107").Value = 1500001920
$ws.Range("L107").Value = 3199.0002
$ws.Range("M107").Value = -1500000000
$ws.Range("N107").Value = -7039.0002
$ws.Range("H132").Value = 1515.9166
$ws.Range("I132").Value = 1079.7188
$ws.Range("J132").Value = 2388.3125
$ws.Range("K132").Value = 3239.1564
$ws.Range("L132").Value = 7164.9375
$ws.Range("M132").Value = -709.1564000000003
$ws.Range("N132").Value = -12224.9375
$ws.Range("H136").Value = 1491.6271
$ws.Range("I136").Value = 863.36365
$ws.Range("J136").Value = 2289.0386
$ws.Range("K136").Value = 2590.09095
$ws.Range("L136").Value = 6867.1158
$ws.Range("M136").Value = -40.09094999999979
$ws.Range("N136").Value = -11967.1158

Write-Host "Applied all changes"